$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "time" column (F) is removed entirely; datetime/empty_column/extra-data
# columns all shift one column to the left.
$ws.Columns("F").Delete() | Out-Null

# Leave the selection on the column that used to hold "time" (now "datetime"),
# matching the workbook's last-saved view state.
$ws.Columns("F").Select() | Out-Null
